$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("B2").Value = 0.4349759063693696
$ws.Range("C2").Value = 0.6084958318088399
$ws.Range("D2").Value = 0.6948499821685598
$ws.Range("E2").Value = 0.7483408004336237
$ws.Range("B3").Value = 0.4495895168240502
$ws.Range("C3").Value = 0.6202535768303969
$ws.Range("D3").Value = 0.7010362124351743
$ws.Range("E3").Value = 0.7522648609778715
$ws.Range("B4").Value = 0.4126039193304918
$ws.Range("C4").Value = 0.5920925721075019
$ws.Range("D4").Value = 0.6867376105788219
$ws.Range("E4").Value = 0.7427419605067365
$ws.Range("B5").Value = 0.4740925553094573
$ws.Range("C5").Value = 0.6404928724557114
$ws.Range("D5").Value = 0.7046736107485284
$ws.Range("E5").Value = 0.7596649851233738
$ws.Range("B6").Value = 0.4692799780971428
$ws.Range("C6").Value = 0.6363875332568553
$ws.Range("D6").Value = 0.7027173425921137
$ws.Range("E6").Value = 0.7583094459532475
$ws.Range("B7").Value = 0.481179764654688
$ws.Range("C7").Value = 0.644483405237085
$ws.Range("D7").Value = 0.7130282815819784
$ws.Range("E7").Value = 0.7601619967263648
$ws.Range("B8").Value = 0.4171685643628851
$ws.Range("C8").Value = 0.5876864638500461
$ws.Range("D8").Value = 0.6848147166973333
$ws.Range("E8").Value = 0.7391238751082457
$ws.Range("B9").Value = 0.4756995253731762
$ws.Range("C9").Value = 0.63951235732726
$ws.Range("D9").Value = 0.7103386006351842
$ws.Range("E9").Value = 0.7585389402505501
$ws.Range("B10").Value = 0.4850680614742685
$ws.Range("C10").Value = 0.645054067969026
$ws.Range("D10").Value = 0.7119964603128214
$ws.Range("E10").Value = 0.7585068958190951
$ws.Range("B11").Value = 0.4846872260551944
$ws.Range("C11").Value = 0.6454288995183227
$ws.Range("D11").Value = 0.7126442726389244
$ws.Range("E11").Value = 0.7591333756474891
$ws.Range("B12").Value = 0.4751725270085909
$ws.Range("C12").Value = 0.6358278061312898
$ws.Range("D12").Value = 0.7063020607180248
$ws.Range("E12").Value = 0.7528757606733335
$ws.Range("B13").Value = 0.4840699376551185
$ws.Range("C13").Value = 0.6439606504256168
$ws.Range("D13").Value = 0.7111449616502787
$ws.Range("E13").Value = 0.7577073778010187
